$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 248
$ws.Range("F3").Value = 856
$ws.Range("F5").Value = 2290
$ws.Range("F6").Value = 1363
$ws.Range("F8").Value = 810
$ws.Range("F9").Value = 1158
$ws.Range("F10").Value = 1043
$ws.Range("F11").Value = 3026
$ws.Range("F15").Value = 617
$ws.Range("F17").Value = 239
$ws.Range("F19").Value = 1094
$ws.Range("F20").Value = 1094
$ws.Range("F21").Value = 149
$ws.Range("F22").Value = 534
$ws.Range("F23").Value = 180
$ws.Range("F25").Value = 230
$ws.Range("F26").Value = 643
$ws.Range("F27").Value = 605
$ws.Range("F29").Value = 838
$ws.Range("F30").Value = 72
$ws.Range("F32").Value = 35
$ws.Range("F33").Value = 1043
$ws.Range("F34").Value = 5046
$ws.Range("F35").Value = 511
$ws.Range("F36").Value = 245
$ws.Range("F37").Value = 133
$ws.Range("F39").Value = 4
$ws.Range("F40").Value = 2

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 363
$ws.Range("F21").Value = 41
$ws.Range("F25").Value = 377
$ws.Range("F28").Value = 678
$ws.Range("F31").Value = 9
$ws.Range("F37").Value = 439
$ws.Range("F39").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 641
$ws.Range("F5").Value = 430
$ws.Range("F6").Value = 411

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 641
$ws.Range("F3").Value = 248
$ws.Range("F4").Value = 430
$ws.Range("F5").Value = 856
$ws.Range("F8").Value = 2290
$ws.Range("F9").Value = 1363
$ws.Range("F11").Value = 810
$ws.Range("F12").Value = 1158
$ws.Range("F15").Value = 1043
$ws.Range("F16").Value = 3026
$ws.Range("F20").Value = 617
$ws.Range("F21").Value = 411
$ws.Range("F23").Value = 239
$ws.Range("F25").Value = 1094
$ws.Range("F26").Value = 1094
$ws.Range("F27").Value = 149
$ws.Range("F28").Value = 41
$ws.Range("F29").Value = 534
$ws.Range("F31").Value = 180
$ws.Range("F32").Value = 230
$ws.Range("F34").Value = 643
$ws.Range("F35").Value = 605
$ws.Range("F36").Value = 377
$ws.Range("F37").Value = 678
$ws.Range("F38").Value = 839
$ws.Range("F40").Value = 72
$ws.Range("F42").Value = 9
$ws.Range("F43").Value = 1043
$ws.Range("F44").Value = 5046
$ws.Range("F46").Value = 511
$ws.Range("F47").Value = 439
$ws.Range("F48").Value = 439
$ws.Range("F49").Value = 245

